$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5518
$wsExhibit.Range("F7").Value = 32
$wsExhibit.Range("F8").Value = 360

# Sheet "演出" (performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 47

# Sheet "全部类型" (all types, aggregated view)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5518
$wsAll.Range("F7").Value = 32
$wsAll.Range("F8").Value = 47
$wsAll.Range("F9").Value = 360
